$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4938.4443
$ws.Range("I86").Value = 5748
$ws.Range("K86").Value = 5748
$ws.Range("M86").Value = -4625
$ws.Range("H89").Value = 4938.4443
$ws.Range("I89").Value = 5748
$ws.Range("K89").Value = 28740
$ws.Range("M89").Value = -23124
$ws.Range("H107").Value = 185.93333
$ws.Range("I107").Value = 185.93333
$ws.Range("K107").Value = 185.93333
$ws.Range("M107").Value = 1734.06667
$ws.Range("H113").Value = 6181.8184
$ws.Range("I113").Value = 4842
$ws.Range("K113").Value = 4842
$ws.Range("M113").Value = -1588
$ws.Range("H125").Value = 9570.333000000001
$ws.Range("J125").Value = 15473.714
$ws.Range("L125").Value = 139263.426
$ws.Range("N125").Value = -144183.426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2944.2
$ws.Range("I61").Value = 1554.619
$ws.Range("J61").Value = 6186.5557
$ws.Range("K61").Value = 1554.619
$ws.Range("L61").Value = 6186.5557
$ws.Range("M61").Value = -1342.619
$ws.Range("N61").Value = -6610.5557
$ws.Range("H88").Value = 1158.95
$ws.Range("I88").Value = 1266.8572
$ws.Range("J88").Value = 1100.8462
$ws.Range("K88").Value = 1266.8572
$ws.Range("L88").Value = 1100.8462
$ws.Range("M88").Value = -860.8571999999999
$ws.Range("N88").Value = -1912.8462
$ws.Range("H91").Value = 1158.95
$ws.Range("I91").Value = 1266.8572
$ws.Range("J91").Value = 1100.8462
$ws.Range("K91").Value = 1266.8572
$ws.Range("L91").Value = 1100.8462
$ws.Range("M91").Value = 137.1428000000001
$ws.Range("N91").Value = -3908.8462
$ws.Range("H136").Value = 2944.2
$ws.Range("I136").Value = 1554.619
$ws.Range("J136").Value = 6186.5557
$ws.Range("K136").Value = 4663.857
$ws.Range("L136").Value = 18559.6671
$ws.Range("M136").Value = -2113.857
$ws.Range("N136").Value = -23659.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3632.9443
$ws.Range("I94").Value = 3024.5625
$ws.Range("K94").Value = 3024.5625
$ws.Range("M94").Value = -2573.5625
$ws.Range("H134").Value = 3614
$ws.Range("I134").Value = 2728.7407
$ws.Range("K134").Value = 8186.222099999999
$ws.Range("M134").Value = -5651.222099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3281.95
$ws.Range("I31").Value = 1290.4
$ws.Range("K31").Value = 1290.4
$ws.Range("M31").Value = -995.4000000000001
$ws.Range("H34").Value = 3281.95
$ws.Range("I34").Value = 1290.4
$ws.Range("K34").Value = 1290.4
$ws.Range("M34").Value = -1088.4
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H122").Value = 643597.75
$ws.Range("I122").Value = 2042715.8
$ws.Range("J122").Value = 7635
$ws.Range("K122").Value = 6128147.4
$ws.Range("L122").Value = 22905
$ws.Range("M122").Value = -6125697.4
$ws.Range("N122").Value = -27805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 10743.947
$ws.Range("I5").Value = 4980.5454
$ws.Range("J5").Value = 18668.625
$ws.Range("K5").Value = 14941.6362
$ws.Range("L5").Value = 56005.875
$ws.Range("M5").Value = -14829.6362
$ws.Range("N5").Value = -56229.875
$ws.Range("H107").Value = 524.75
$ws.Range("J107").Value = 524.75
$ws.Range("L107").Value = 1574.25
$ws.Range("N107").Value = -5414.25
$ws.Range("H121").Value = 2353.9375
$ws.Range("J121").Value = 2490.8667
$ws.Range("L121").Value = 7472.6001
$ws.Range("N121").Value = -10092.6001
$ws.Range("H135").Value = 10743.947
$ws.Range("I135").Value = 4980.5454
$ws.Range("J135").Value = 18668.625
$ws.Range("K135").Value = 44824.9086
$ws.Range("L135").Value = 168017.625
$ws.Range("M135").Value = -42289.9086
$ws.Range("N135").Value = -173087.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 69073.47
$ws.Range("I70").Value = 104105
$ws.Range("J70").Value = 4849
$ws.Range("K70").Value = 104105
$ws.Range("L70").Value = 4849
$ws.Range("M70").Value = -103835
$ws.Range("N70").Value = -5389
$ws.Range("H73").Value = 69073.47
$ws.Range("I73").Value = 104105
$ws.Range("J73").Value = 4849
$ws.Range("K73").Value = 104105
$ws.Range("L73").Value = 4849
$ws.Range("M73").Value = -103169
$ws.Range("N73").Value = -6721
$ws.Range("H102").Value = 1112.4138
$ws.Range("I102").Value = 827.7778
$ws.Range("K102").Value = 827.7778
$ws.Range("M102").Value = 794.2222
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -36884
$ws.Range("H132").Value = 3341.0833
$ws.Range("I132").Value = 1982.6875
$ws.Range("K132").Value = 5948.0625
$ws.Range("M132").Value = -3418.0625
$ws.Range("H138").Value = 69922.75
$ws.Range("J138").Value = 69922.75
$ws.Range("L138").Value = 69922.75
$ws.Range("N138").Value = -80202.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4055.4285
$ws.Range("I7").Value = 2610.7273
$ws.Range("K7").Value = 2610.7273
$ws.Range("M7").Value = -2498.7273
$ws.Range("H46").Value = 2646.875
$ws.Range("I46").Value = 1115
$ws.Range("K46").Value = 1115
$ws.Range("M46").Value = -927
$ws.Range("H61").Value = 2845.468
$ws.Range("I61").Value = 2661.353
$ws.Range("K61").Value = 2661.353
$ws.Range("M61").Value = -2459.353
$ws.Range("H93").Value = 1430874.1
$ws.Range("I93").Value = 2686.6667
$ws.Range("J93").Value = 9999999
$ws.Range("K93").Value = 2686.6667
$ws.Range("L93").Value = 9999999
$ws.Range("M93").Value = -1438.6667
$ws.Range("N93").Value = -10002495
$ws.Range("H113").Value = 2845.468
$ws.Range("I113").Value = 2661.353
$ws.Range("K113").Value = 2661.353
$ws.Range("M113").Value = -491.3530000000001
$ws.Range("H118").Value = 43615
$ws.Range("J118").Value = 43615
$ws.Range("L118").Value = 43615
$ws.Range("N118").Value = -46929
$ws.Range("H122").Value = 5225.8
$ws.Range("I122").Value = 2784.6667
$ws.Range("K122").Value = 8354.000100000001
$ws.Range("M122").Value = -5904.000100000001
$ws.Range("H126").Value = 4055.4285
$ws.Range("I126").Value = 2610.7273
$ws.Range("K126").Value = 7832.1819
$ws.Range("M126").Value = -5362.1819
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2020.95
$ws.Range("I126").Value = 1770.5
$ws.Range("K126").Value = 5311.5
$ws.Range("M126").Value = -2841.5
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H136").Value = 2406.476
$ws.Range("I136").Value = 1399.1936
$ws.Range("K136").Value = 4197.5808
$ws.Range("M136").Value = -1647.5808
